$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lists")
$ws.Activate()

# Shift "Others " down from B24 to B26, and insert the two new truck
# types "Car Shuttler" (B24) and "Lowbed" (B25) above it.
$ws.Range("B26").Value = "Others "
$ws.Range("B25").Value = "Lowbed"
$ws.Range("B24").Value = "Car Shuttler"

# Copy B24's formatting onto the two newly-populated cells so they pick
# up the same style as the rest of the list.
$ws.Range("B24").Copy()
$ws.Range("B25:B26").PasteSpecial(-4122)

# Grow the TruckTypes named range to cover the two new rows.
$wb.Names.Item("TruckTypes").RefersTo = "=Lists!`$B`$2:`$B`$26"

[void]$ws.Range("B25").Select()

# Return to the main Trucks sheet, leaving it the active tab/selection.
$ws1 = $wb.Worksheets.Item("Trucks")
$ws1.Activate()
[void]$ws1.Range("F10").Select()
